$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.385511999999999
$ws.Range("H2").Value = 16.156536
$ws.Range("I2").Value = 0.02736372477514656
$ws.Range("J2").Value = 0.02736372477514657
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.01650666666666667
$ws.Range("N2").Value = 0.04952
$ws.Range("O2").Value = 0.795859985214233
$ws.Range("P2").Value = 0.795859985214233
$ws.Range("Q2").Value = 0.08889685141333332
$ws.Range("R2").Value = 0.80007166272
$ws.Range("S2").Value = 0.02177769359495449
$ws.Range("T2").Value = 0.02177769359495449

# Row 3
$ws.Range("G3").Value = 5.385511999999999
$ws.Range("H3").Value = 16.156536
$ws.Range("I3").Value = 0.02736372477514656
$ws.Range("J3").Value = 0.02736372477514657
$ws.Range("O3").Value = 0.2041400147857671
$ws.Range("P3").Value = 0.2041400147857671
$ws.Range("Q3").Value = 0.022802257808
$ws.Range("R3").Value = 0.205220320272
$ws.Range("S3").Value = 0.005586031180192081
$ws.Range("T3").Value = 0.005586031180192082

# Row 4
$ws.Range("I4").Value = 0.7812411799860843
$ws.Range("J4").Value = 0.7812411799860843
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.01650666666666667
$ws.Range("N4").Value = 0.04952
$ws.Range("O4").Value = 0.795859985214233
$ws.Range("P4").Value = 0.795859985214233
$ws.Range("Q4").Value = 2.538027321422222
$ws.Range("R4").Value = 22.8422458928
$ws.Range("S4").Value = 0.621758593952475
$ws.Range("T4").Value = 0.621758593952475

# Row 5
$ws.Range("I5").Value = 0.7812411799860843
$ws.Range("J5").Value = 0.7812411799860843
$ws.Range("O5").Value = 0.2041400147857671
$ws.Range("P5").Value = 0.2041400147857671
$ws.Range("S5").Value = 0.1594825860336094
$ws.Range("T5").Value = 0.1594825860336094

# Row 6
$ws.Range("I6").Value = 0.1913950952387691
$ws.Range("J6").Value = 0.1913950952387691
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.01650666666666667
$ws.Range("N6").Value = 0.04952
$ws.Range("O6").Value = 0.795859985214233
$ws.Range("P6").Value = 0.795859985214233
$ws.Range("Q6").Value = 0.6217874752977778
$ws.Range("R6").Value = 5.59608727768
$ws.Range("S6").Value = 0.1523236976668035
$ws.Range("T6").Value = 0.1523236976668035

# Row 7
$ws.Range("I7").Value = 0.1913950952387691
$ws.Range("J7").Value = 0.1913950952387691
$ws.Range("O7").Value = 0.2041400147857671
$ws.Range("P7").Value = 0.2041400147857671
$ws.Range("S7").Value = 0.03907139757196562
$ws.Range("T7").Value = 0.03907139757196562
